# Auto-generated COM-interop script implementing the target diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Fix up the existing 'ODI Batting' sheet (currently sheet index 1):
#    MATCH_CARD_LINK -> MATCH_CODE, URL values -> bare numeric codes.
# ---------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"
$odiBatting.Range("D2").NumberFormat = "@"
$odiBatting.Range("D2").Value = '3149'
$odiBatting.Range("D3").NumberFormat = "@"
$odiBatting.Range("D3").Value = '3153'
$odiBatting.Range("D4").NumberFormat = "@"
$odiBatting.Range("D4").Value = '3164'
$odiBatting.Range("D5").NumberFormat = "@"
$odiBatting.Range("D5").Value = '3166'
$odiBatting.Range("D6").NumberFormat = "@"
$odiBatting.Range("D6").Value = '3183'
$odiBatting.Range("D7").NumberFormat = "@"
$odiBatting.Range("D7").Value = '3187'
$odiBatting.Range("D8").NumberFormat = "@"
$odiBatting.Range("D8").Value = '3406'
$odiBatting.Range("D9").NumberFormat = "@"
$odiBatting.Range("D9").Value = '3407'
$odiBatting.Range("D10").NumberFormat = "@"
$odiBatting.Range("D10").Value = '3421'
$odiBatting.Range("D11").NumberFormat = "@"
$odiBatting.Range("D11").Value = '3437'
$odiBatting.Range("D12").NumberFormat = "@"
$odiBatting.Range("D12").Value = '3651'
$odiBatting.Range("D13").NumberFormat = "@"
$odiBatting.Range("D13").Value = '3652'
$odiBatting.Range("D14").NumberFormat = "@"
$odiBatting.Range("D14").Value = '3705'
$odiBatting.Range("D15").NumberFormat = "@"
$odiBatting.Range("D15").Value = '3709'
$odiBatting.Range("D16").NumberFormat = "@"
$odiBatting.Range("D16").Value = '3711'
$odiBatting.Range("D17").NumberFormat = "@"
$odiBatting.Range("D17").Value = '3721'
$odiBatting.Range("D18").NumberFormat = "@"
$odiBatting.Range("D18").Value = '3722'
$odiBatting.Range("D19").NumberFormat = "@"
$odiBatting.Range("D19").Value = '3725'
$odiBatting.Range("D20").NumberFormat = "@"
$odiBatting.Range("D20").Value = '3730'
$odiBatting.Range("D21").NumberFormat = "@"
$odiBatting.Range("D21").Value = '3754'
$odiBatting.Range("D22").NumberFormat = "@"
$odiBatting.Range("D22").Value = '3759'
$odiBatting.Range("D23").NumberFormat = "@"
$odiBatting.Range("D23").Value = '3764'
$odiBatting.Range("D24").NumberFormat = "@"
$odiBatting.Range("D24").Value = '3773'
$odiBatting.Range("D25").NumberFormat = "@"
$odiBatting.Range("D25").Value = '3778'
$odiBatting.Range("D26").NumberFormat = "@"
$odiBatting.Range("D26").Value = '3785'
$odiBatting.Range("D27").NumberFormat = "@"
$odiBatting.Range("D27").Value = '3912'
$odiBatting.Range("D28").NumberFormat = "@"
$odiBatting.Range("D28").Value = '3913'
$odiBatting.Range("D29").NumberFormat = "@"
$odiBatting.Range("D29").Value = '3915'
$odiBatting.Range("D30").NumberFormat = "@"
$odiBatting.Range("D30").Value = '3916'
$odiBatting.Range("D31").NumberFormat = "@"
$odiBatting.Range("D31").Value = '3918'
$odiBatting.Range("D32").NumberFormat = "@"
$odiBatting.Range("D32").Value = '4040'
$odiBatting.Range("D33").NumberFormat = "@"
$odiBatting.Range("D33").Value = '4043'
$odiBatting.Range("D34").NumberFormat = "@"
$odiBatting.Range("D34").Value = '4046'
$odiBatting.Range("D35").NumberFormat = "@"
$odiBatting.Range("D35").Value = '4092'
$odiBatting.Range("D36").NumberFormat = "@"
$odiBatting.Range("D36").Value = '4093'
$odiBatting.Range("D37").NumberFormat = "@"
$odiBatting.Range("D37").Value = '4097'
$odiBatting.Range("D38").NumberFormat = "@"
$odiBatting.Range("D38").Value = '4136'
$odiBatting.Range("D39").NumberFormat = "@"
$odiBatting.Range("D39").Value = '4147'
$odiBatting.Range("D40").NumberFormat = "@"
$odiBatting.Range("D40").Value = '4154'
$odiBatting.Range("D41").NumberFormat = "@"
$odiBatting.Range("D41").Value = '4203'
$odiBatting.Range("D42").NumberFormat = "@"
$odiBatting.Range("D42").Value = '4259'
$odiBatting.Range("D43").NumberFormat = "@"
$odiBatting.Range("D43").Value = '4267'
$odiBatting.Range("D44").NumberFormat = "@"
$odiBatting.Range("D44").Value = '4377'
$odiBatting.Range("D45").NumberFormat = "@"
$odiBatting.Range("D45").Value = '4378'
$odiBatting.Range("D46").NumberFormat = "@"
$odiBatting.Range("D46").Value = '4444'
$odiBatting.Range("D47").NumberFormat = "@"
$odiBatting.Range("D47").Value = '4446'
$odiBatting.Range("D48").NumberFormat = "@"
$odiBatting.Range("D48").Value = '4448'

# ---------------------------------------------------------------
# 2. Fix up the existing 'ODI Bowling' sheet (currently sheet index 2):
#    MATCH_CARD_LINK -> MATCH_CODE, URL values -> bare numeric codes.
# ---------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"
$odiBowling.Range("B2").NumberFormat = "@"
$odiBowling.Range("B2").Value = '3187'
$odiBowling.Range("B3").NumberFormat = "@"
$odiBowling.Range("B3").Value = '3407'
$odiBowling.Range("B4").NumberFormat = "@"
$odiBowling.Range("B4").Value = '3421'
$odiBowling.Range("B5").NumberFormat = "@"
$odiBowling.Range("B5").Value = '3437'
$odiBowling.Range("B6").NumberFormat = "@"
$odiBowling.Range("B6").Value = '3705'
$odiBowling.Range("B7").NumberFormat = "@"
$odiBowling.Range("B7").Value = '3709'
$odiBowling.Range("B8").NumberFormat = "@"
$odiBowling.Range("B8").Value = '3722'
$odiBowling.Range("B9").NumberFormat = "@"
$odiBowling.Range("B9").Value = '3730'
$odiBowling.Range("B10").NumberFormat = "@"
$odiBowling.Range("B10").Value = '3754'
$odiBowling.Range("B11").NumberFormat = "@"
$odiBowling.Range("B11").Value = '3764'
$odiBowling.Range("B12").NumberFormat = "@"
$odiBowling.Range("B12").Value = '3773'
$odiBowling.Range("B13").NumberFormat = "@"
$odiBowling.Range("B13").Value = '3913'
$odiBowling.Range("B14").NumberFormat = "@"
$odiBowling.Range("B14").Value = '3915'
$odiBowling.Range("B15").NumberFormat = "@"
$odiBowling.Range("B15").Value = '4092'
$odiBowling.Range("B16").NumberFormat = "@"
$odiBowling.Range("B16").Value = '4203'
$odiBowling.Range("B17").NumberFormat = "@"
$odiBowling.Range("B17").Value = '4377'
$odiBowling.Range("B18").NumberFormat = "@"
$odiBowling.Range("B18").Value = '4378'
$odiBowling.Range("B19").NumberFormat = "@"
$odiBowling.Range("B19").Value = '4446'
$odiBowling.Range("B20").NumberFormat = "@"
$odiBowling.Range("B20").Value = '4448'

# ---------------------------------------------------------------
# 3. Insert the new 'Player Info' sheet as the first tab.
# ---------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

# ---------------------------------------------------------------
# 4. Insert the new 'ODI Batting Extra' sheet as the last tab.
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

# ---------------------------------------------------------------
# 5. Populate 'Player Info' (headers get the same bold / thin-border /
#    center-top-aligned look as the existing header row).
# ---------------------------------------------------------------
$playerInfo.Range("A1").Value = 'ID'
$playerInfo.Range("B1").Value = 'NAME'
$playerInfo.Range("C1").Value = 'BATTING_HAND'
$playerInfo.Range("D1").Value = 'BOWL_STYLE'
$piHeader = $playerInfo.Range("A1:D1")
$piHeader.Font.Bold = $true
$piHeader.HorizontalAlignment = -4108
$piHeader.VerticalAlignment = -4160
$piHeader.Borders.LineStyle = 1
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = '3813'
$playerInfo.Range("B2").NumberFormat = "@"
$playerInfo.Range("B2").Value = 'Javed Ahmadi'
$playerInfo.Range("C2").NumberFormat = "@"
$playerInfo.Range("C2").Value = 'Right Handed'
$playerInfo.Range("D2").NumberFormat = "@"
$playerInfo.Range("D2").Value = 'Right Arm Off Break'

# ---------------------------------------------------------------
# 6. Populate 'ODI Batting Extra'.
# ---------------------------------------------------------------
$battingExtra.Range("A1").Value = 'MATCH_CODE'
$battingExtra.Range("B1").Value = 'BATTING_POSITION'
$battingExtra.Range("C1").Value = 'NUM_4'
$battingExtra.Range("D1").Value = 'NUM_6'
$battingExtra.Range("E1").Value = 'PERCENT_RUNS_OF_TOTAL'
$battingExtra.Range("F1").Value = 'MAN_OF_MATCH'
$beHeader = $battingExtra.Range("A1:F1")
$beHeader.Font.Bold = $true
$beHeader.HorizontalAlignment = -4108
$beHeader.VerticalAlignment = -4160
$beHeader.Borders.LineStyle = 1
$battingExtra.Range("A2").NumberFormat = "@"
$battingExtra.Range("A2").Value = '3915'
$battingExtra.Range("B2").Value = 2
$battingExtra.Range("C2").NumberFormat = "@"
$battingExtra.Range("C2").Value = '3'
$battingExtra.Range("D2").NumberFormat = "@"
$battingExtra.Range("D2").Value = '0'
$battingExtra.Range("E2").NumberFormat = "@"
$battingExtra.Range("E2").Value = '11.60%'
$battingExtra.Range("F2").NumberFormat = "@"
$battingExtra.Range("F2").Value = 'NO'
$battingExtra.Range("A3").NumberFormat = "@"
$battingExtra.Range("A3").Value = '3916'
$battingExtra.Range("B3").Value = ""
$battingExtra.Range("C3").Value = ""
$battingExtra.Range("D3").Value = ""
$battingExtra.Range("E3").Value = ""
$battingExtra.Range("F3").NumberFormat = "@"
$battingExtra.Range("F3").Value = 'NO'
$battingExtra.Range("A4").NumberFormat = "@"
$battingExtra.Range("A4").Value = '3918'
$battingExtra.Range("B4").Value = 2
$battingExtra.Range("C4").NumberFormat = "@"
$battingExtra.Range("C4").Value = '0'
$battingExtra.Range("D4").NumberFormat = "@"
$battingExtra.Range("D4").Value = '0'
$battingExtra.Range("E4").NumberFormat = "@"
$battingExtra.Range("E4").Value = '1.19%'
$battingExtra.Range("F4").NumberFormat = "@"
$battingExtra.Range("F4").Value = 'NO'
$battingExtra.Range("A5").NumberFormat = "@"
$battingExtra.Range("A5").Value = '4040'
$battingExtra.Range("B5").Value = 2
$battingExtra.Range("C5").NumberFormat = "@"
$battingExtra.Range("C5").Value = '8'
$battingExtra.Range("D5").NumberFormat = "@"
$battingExtra.Range("D5").Value = '2'
$battingExtra.Range("E5").NumberFormat = "@"
$battingExtra.Range("E5").Value = '38.21%'
$battingExtra.Range("F5").NumberFormat = "@"
$battingExtra.Range("F5").Value = 'NO'
$battingExtra.Range("A6").NumberFormat = "@"
$battingExtra.Range("A6").Value = '4043'
$battingExtra.Range("B6").Value = 2
$battingExtra.Range("C6").NumberFormat = "@"
$battingExtra.Range("C6").Value = '1'
$battingExtra.Range("D6").NumberFormat = "@"
$battingExtra.Range("D6").Value = '0'
$battingExtra.Range("E6").NumberFormat = "@"
$battingExtra.Range("E6").Value = '6.67%'
$battingExtra.Range("F6").NumberFormat = "@"
$battingExtra.Range("F6").Value = 'NO'
$battingExtra.Range("A7").NumberFormat = "@"
$battingExtra.Range("A7").Value = '4046'
$battingExtra.Range("B7").Value = ""
$battingExtra.Range("C7").Value = ""
$battingExtra.Range("D7").Value = ""
$battingExtra.Range("E7").Value = ""
$battingExtra.Range("F7").NumberFormat = "@"
$battingExtra.Range("F7").Value = 'NO'
$battingExtra.Range("A8").NumberFormat = "@"
$battingExtra.Range("A8").Value = '4092'
$battingExtra.Range("B8").Value = 1
$battingExtra.Range("C8").NumberFormat = "@"
$battingExtra.Range("C8").Value = '1'
$battingExtra.Range("D8").NumberFormat = "@"
$battingExtra.Range("D8").Value = '0'
$battingExtra.Range("E8").NumberFormat = "@"
$battingExtra.Range("E8").Value = '4.20%'
$battingExtra.Range("F8").NumberFormat = "@"
$battingExtra.Range("F8").Value = 'NO'
$battingExtra.Range("A9").NumberFormat = "@"
$battingExtra.Range("A9").Value = '4093'
$battingExtra.Range("B9").Value = 1
$battingExtra.Range("C9").NumberFormat = "@"
$battingExtra.Range("C9").Value = '7'
$battingExtra.Range("D9").NumberFormat = "@"
$battingExtra.Range("D9").Value = '0'
$battingExtra.Range("E9").NumberFormat = "@"
$battingExtra.Range("E9").Value = '21.82%'
$battingExtra.Range("F9").NumberFormat = "@"
$battingExtra.Range("F9").Value = 'NO'
$battingExtra.Range("A10").NumberFormat = "@"
$battingExtra.Range("A10").Value = '4097'
$battingExtra.Range("B10").Value = 1
$battingExtra.Range("C10").NumberFormat = "@"
$battingExtra.Range("C10").Value = '3'
$battingExtra.Range("D10").NumberFormat = "@"
$battingExtra.Range("D10").Value = '1'
$battingExtra.Range("E10").NumberFormat = "@"
$battingExtra.Range("E10").Value = '15.25%'
$battingExtra.Range("F10").NumberFormat = "@"
$battingExtra.Range("F10").Value = 'NO'
$battingExtra.Range("A11").NumberFormat = "@"
$battingExtra.Range("A11").Value = '4136'
$battingExtra.Range("B11").Value = 2
$battingExtra.Range("C11").NumberFormat = "@"
$battingExtra.Range("C11").Value = '6'
$battingExtra.Range("D11").NumberFormat = "@"
$battingExtra.Range("D11").Value = '2'
$battingExtra.Range("E11").NumberFormat = "@"
$battingExtra.Range("E11").Value = '31.54%'
$battingExtra.Range("F11").NumberFormat = "@"
$battingExtra.Range("F11").Value = 'NO'
$battingExtra.Range("A12").NumberFormat = "@"
$battingExtra.Range("A12").Value = '4147'
$battingExtra.Range("B12").Value = ""
$battingExtra.Range("C12").Value = ""
$battingExtra.Range("D12").Value = ""
$battingExtra.Range("E12").Value = ""
$battingExtra.Range("F12").NumberFormat = "@"
$battingExtra.Range("F12").Value = 'NO'
$battingExtra.Range("A13").NumberFormat = "@"
$battingExtra.Range("A13").Value = '4154'
$battingExtra.Range("B13").Value = ""
$battingExtra.Range("C13").Value = ""
$battingExtra.Range("D13").Value = ""
$battingExtra.Range("E13").Value = ""
$battingExtra.Range("F13").NumberFormat = "@"
$battingExtra.Range("F13").Value = 'NO'
$battingExtra.Range("A14").NumberFormat = "@"
$battingExtra.Range("A14").Value = '4203'
$battingExtra.Range("B14").Value = ""
$battingExtra.Range("C14").Value = ""
$battingExtra.Range("D14").Value = ""
$battingExtra.Range("E14").Value = ""
$battingExtra.Range("F14").NumberFormat = "@"
$battingExtra.Range("F14").Value = 'NO'
$battingExtra.Range("A15").NumberFormat = "@"
$battingExtra.Range("A15").Value = '4259'
$battingExtra.Range("B15").Value = ""
$battingExtra.Range("C15").Value = ""
$battingExtra.Range("D15").Value = ""
$battingExtra.Range("E15").Value = ""
$battingExtra.Range("F15").NumberFormat = "@"
$battingExtra.Range("F15").Value = 'NO'
$battingExtra.Range("A16").NumberFormat = "@"
$battingExtra.Range("A16").Value = '4267'
$battingExtra.Range("B16").Value = 2
$battingExtra.Range("C16").NumberFormat = "@"
$battingExtra.Range("C16").Value = '2'
$battingExtra.Range("D16").NumberFormat = "@"
$battingExtra.Range("D16").Value = '1'
$battingExtra.Range("E16").NumberFormat = "@"
$battingExtra.Range("E16").Value = '11.11%'
$battingExtra.Range("F16").NumberFormat = "@"
$battingExtra.Range("F16").Value = 'NO'
$battingExtra.Range("A17").NumberFormat = "@"
$battingExtra.Range("A17").Value = '4377'
$battingExtra.Range("B17").Value = 2
$battingExtra.Range("C17").NumberFormat = "@"
$battingExtra.Range("C17").Value = '1'
$battingExtra.Range("D17").NumberFormat = "@"
$battingExtra.Range("D17").Value = '0'
$battingExtra.Range("E17").NumberFormat = "@"
$battingExtra.Range("E17").Value = '2.58%'
$battingExtra.Range("F17").NumberFormat = "@"
$battingExtra.Range("F17").Value = 'NO'
$battingExtra.Range("A18").NumberFormat = "@"
$battingExtra.Range("A18").Value = '4378'
$battingExtra.Range("B18").Value = 2
$battingExtra.Range("C18").NumberFormat = "@"
$battingExtra.Range("C18").Value = '0'
$battingExtra.Range("D18").NumberFormat = "@"
$battingExtra.Range("D18").Value = '0'
$battingExtra.Range("E18").Value = ""
$battingExtra.Range("F18").NumberFormat = "@"
$battingExtra.Range("F18").Value = 'NO'
$battingExtra.Range("A19").NumberFormat = "@"
$battingExtra.Range("A19").Value = '4444'
$battingExtra.Range("B19").Value = 2
$battingExtra.Range("C19").NumberFormat = "@"
$battingExtra.Range("C19").Value = '2'
$battingExtra.Range("D19").NumberFormat = "@"
$battingExtra.Range("D19").Value = '2'
$battingExtra.Range("E19").NumberFormat = "@"
$battingExtra.Range("E19").Value = '13.24%'
$battingExtra.Range("F19").NumberFormat = "@"
$battingExtra.Range("F19").Value = 'NO'
$battingExtra.Range("A20").NumberFormat = "@"
$battingExtra.Range("A20").Value = '4446'
$battingExtra.Range("B20").Value = 2
$battingExtra.Range("C20").NumberFormat = "@"
$battingExtra.Range("C20").Value = '2'
$battingExtra.Range("D20").NumberFormat = "@"
$battingExtra.Range("D20").Value = '0'
$battingExtra.Range("E20").NumberFormat = "@"
$battingExtra.Range("E20").Value = '6.15%'
$battingExtra.Range("F20").NumberFormat = "@"
$battingExtra.Range("F20").Value = 'NO'
$battingExtra.Range("A21").NumberFormat = "@"
$battingExtra.Range("A21").Value = '4448'
$battingExtra.Range("B21").Value = 2
$battingExtra.Range("C21").NumberFormat = "@"
$battingExtra.Range("C21").Value = '3'
$battingExtra.Range("D21").NumberFormat = "@"
$battingExtra.Range("D21").Value = '1'
$battingExtra.Range("E21").NumberFormat = "@"
$battingExtra.Range("E21").Value = '7.14%'
$battingExtra.Range("F21").NumberFormat = "@"
$battingExtra.Range("F21").Value = 'NO'
